$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Vida
$ws.Range("D4").Value = 2000
$ws.Range("F4").Value = 2000
$ws.Range("H4").Value = 3000
$ws.Range("J4").Value = 3000
$ws.Range("L4").Value = 4000
$ws.Range("N4").Value = 4000

# Row 5 - Ataque
$ws.Range("D5").Value = 500
$ws.Range("F5").Value = 500
$ws.Range("H5").Value = 600
$ws.Range("J5").Value = 600
$ws.Range("L5").Value = 700
$ws.Range("N5").Value = 700

# Row 6 - Defesa
$ws.Range("D6").Value = 300
$ws.Range("F6").Value = 300
$ws.Range("H6").Value = 200
$ws.Range("J6").Value = 200
$ws.Range("L6").Value = 200
$ws.Range("N6").Value = 200

# Row 7 - Mana (power moved from Clérigos/Druidas to Mago/Feiticeiros)
$ws.Range("H7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("L7").Value = 800
$ws.Range("N7").Value = 800

# Row 8 - Fé (power moved from Mago/Feiticeiros to Clérigos/Druidas)
$ws.Range("H8").Value = 800
$ws.Range("J8").Value = 800
$ws.Range("L8").Value = ""
$ws.Range("N8").Value = ""

# Row 13 - Custo de Mana
$ws.Range("D13").Value = 200
$ws.Range("F13").Value = 200

# Row 23 - Custo de Fe
$ws.Range("D23").Value = 500
$ws.Range("F23").Value = 400

# Selection moved from Q10 to H22, scrolled back to top-left (A1 area)
$ws.Range("H22").Select()
